$d = $word.ActiveDocument

# 1. Body text: bump version "These are known issues in Lightning 1.1.1." -> "...1.1.2."
$body = $d.Content
$body.Find.Execute("Lightning 1.1.1.", $false, $false, $false, $false, $false, `
                    $true, 1, $false, "Lightning 1.1.2.", 2)

# 2. Header title: bump version "Lightning Known Issues for 1.1.1" -> "...1.1.2"
$header = $d.Sections(1).Headers(1).Range
$header.Find.Execute("Known Issues for 1.1.1", $false, $false, $false, $false, $false, `
                      $true, 1, $false, "Known Issues for 1.1.2", 2)

# 3. Header date: "October 30, 2022" -> "November 5, 2022"
$header2 = $d.Sections(1).Headers(1).Range
$header2.Find.Execute("October 30, 2022", $false, $false, $false, $false, $false, `
                       $true, 1, $false, "November 5, 2022", 2)
